# "rows met no changes verwijderd" - remove the "NO CHANGE" markers in
# column K (TypeChange) for every row that had no change. Rows 9 and 10
# (which have real TypeChange values) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("translations")

$rows = @(2, 3, 4, 5, 6, 7, 8, 11, 12, 13, 14)
foreach ($r in $rows) {
    $ws.Range("K$r").ClearContents()
}
